$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V2").Value  = 43.96
$ws.Range("V3").Value  = 45.59
$ws.Range("V4").Value  = 49.36
$ws.Range("V5").Value  = 62.07
$ws.Range("V6").Value  = 63.76
$ws.Range("V7").Value  = 18.49
$ws.Range("V8").Value  = 42.06
$ws.Range("V9").Value  = 43.56
$ws.Range("V10").Value = 55.66
$ws.Range("V11").Value = 67.66
$ws.Range("V12").Value = 83.2
$ws.Range("V15").Value = 57.26
$ws.Range("V18").Value = 56.69
$ws.Range("V19").Value = 56.86
$ws.Range("V20").Value = 58.26
$ws.Range("V21").Value = 54.05
